$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date-serial "last changed" timestamp for
# every data row (starting at row 2). The whole column was bumped forward
# by exactly one day (e.g. 45178 -> 45179) when the source data refreshed.

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($lastRow -lt 2) {
    $lastRow = 1
}

if ($lastRow -ge 2) {
    $rng = $ws.Range("C2:C$lastRow")
    $vals = $rng.Value2

    if ($rng.Rows.Count -eq 1) {
        if ($vals -ne $null) {
            $rng.Value2 = $vals + 1
        }
    } else {
        for ($i = 1; $i -le $vals.GetLength(0); $i++) {
            if ($vals[$i, 1] -ne $null) {
                $vals[$i, 1] = $vals[$i, 1] + 1
            }
        }
        $rng.Value2 = $vals
    }
}
